$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Formula = '="' + $value.Replace('"', '""') + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

Set-TextValue 2 4 '25.756.62'
Set-TextValue 2 5 '  -0.40%  '
Set-TextValue 3 4 '1.581.98'
Set-TextValue 3 5 '  -2.35%  '
Set-TextValue 4 4 '1.00'
Set-TextValue 4 5 '  -0.49%  '
Set-TextValue 5 4 '208.82'
Set-TextValue 5 5 '  -1.88%  '
Set-TextValue 6 4 '1.00'
Set-TextValue 6 5 '  -0.50%  '
Set-TextValue 7 4 '0.482'
Set-TextValue 7 5 '  -3.66%  '
Set-TextValue 8 5 '  -0.78%  '
Set-TextValue 9 4 '0.0617'
Set-TextValue 9 5 '  +0.10%  '
Set-TextValue 10 4 '18.16'
Set-TextValue 10 5 '  -1.67%  '
Set-TextValue 11 4 '0.0786'
Set-TextValue 11 5 '  -0.46%  '
Set-TextValue 12 4 '1.803.06'
Set-TextValue 12 5 '  -2.25%  '
Set-TextValue 13 4 '1.573.23'
Set-TextValue 13 5 '  -3.03%  '
Set-TextValue 14 4 '4.02'
Set-TextValue 14 5 '  -2.81%  '
Set-TextValue 15 4 '0.513'
Set-TextValue 15 5 '  -2.26%  '
Set-TextValue 16 4 '25.765.81'
Set-TextValue 16 5 '  -0.42%  '
Set-TextValue 17 4 '60.15'
Set-TextValue 17 5 '  -2.27%  '
Set-TextValue 18 5 '  -1.82%  '
Set-TextValue 19 4 '0.999'
Set-TextValue 19 5 '  -0.55%  '
Set-TextValue 20 4 '191.08'
Set-TextValue 20 5 '  -0.27%  '
Set-TextValue 21 4 '4.19'
Set-TextValue 21 5 '  -1.30%  '
Set-TextValue 22 4 '9.42'
Set-TextValue 22 5 '  -0.79%  '
Set-TextValue 23 4 '5.93'
Set-TextValue 23 5 '  -1.72%  '
Set-TextValue 24 5 '  -2.96%  '
Set-TextValue 25 4 '141.02'
Set-TextValue 25 5 '  -1.97%  '
Set-TextValue 26 4 '1.00'
Set-TextValue 26 5 '  -0.48%  '
Set-TextValue 27 4 '1.71'
Set-TextValue 27 5 '  -1.15%  '
Set-TextValue 28 4 '15.17'
Set-TextValue 28 5 '  -0.22%  '
Set-TextValue 29 4 '6.46'
Set-TextValue 29 5 '  -2.85%  '
Set-TextValue 30 5 '  -5.80%  '
Set-TextValue 31 5 '  -1.10%  '
Set-TextValue 32 5 '  -0.33%  '
Set-TextValue 33 4 '3.02'
Set-TextValue 33 5 '  -2.43%  '
Set-TextValue 34 5 '  +0.77%  '
Set-TextValue 35 5 '  -4.08%  '
Set-TextValue 36 4 '1.102.00'
Set-TextValue 36 5 '  -2.11%  '
Set-TextValue 37 5 '  -0.64%  '
Set-TextValue 38 4 '0.504'
Set-TextValue 38 5 '  -1.33%  '
Set-TextValue 39 2 'MXToken'
Set-TextValue 39 3 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 39 4 '2.33'
Set-TextValue 39 5 '  -2.00%  '
Set-TextValue 40 2 'VeChain'
Set-TextValue 40 3 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 40 4 '0.0151'
Set-TextValue 40 5 '  -1.55%  '
Set-TextValue 41 4 '0.785'
Set-TextValue 41 5 '  -6.53%  '
Set-TextValue 42 5 '  +7.66%  '
Set-TextValue 43 2 'FraxShare'
Set-TextValue 43 3 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 43 4 '5.15'
Set-TextValue 43 5 '  +2.36%  '
Set-TextValue 44 2 'Quant'
Set-TextValue 44 3 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue 44 4 '93.50'
Set-TextValue 44 5 '  -4.83%  '
Set-TextValue 45 4 '1.717.24'
Set-TextValue 45 5 '  -2.19%  '
Set-TextValue 46 5 '  -1.12%  '
Set-TextValue 47 5 '  -0.91%  '
Set-TextValue 48 4 '53.26'
Set-TextValue 48 5 '  -1.42%  '
Set-TextValue 49 5 '  -1.91%  '
Set-TextValue 50 4 '0.406'
Set-TextValue 50 5 '  -1.23%  '
Set-TextValue 51 4 '0.999'
Set-TextValue 51 5 '  -0.51%  '
